$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# 1. Fix ALPHASITIO's balance: 1203626.22 -> 203626.22 (row 2, column C)
$ws.Cells.Item(2, 3).Value = 203626.22

# 2. Insert a new row for ISABEL right after LUIZ (row 6), pushing HEITOR
#    and all following rows down by one.
$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "005624730"
$ws.Cells.Item(7, 2).Value = "ISABEL"
$ws.Cells.Item(7, 3).Value = 5000
